# Append new "lab.mixture.table.*" translation rows to the "Import" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

$rows = @(
    @("lab.mixture.table.name", "Název"),
    @("lab.mixture.table.code", "Kód"),
    @("lab.mixture.table.steep", "Doba zrání"),
    @("lab.mixture.table.pgvg", "PG/VG"),
    @("lab.mixture.table.nicotine", "Obsah nikotinu"),
    @("lab.mixture.table.mixed", "Datum mixu"),
    @("lab.mixture.table.volume", "Obsah"),
    @("lab.mixture.table.liquid", "Liquid"),
    @("lab.mixture.table.booster", "Booster"),
    @("lab.mixture.table.base", "Báze"),
    @("lab.mixture.table.expires", "Datum expirace")
)

$startRow = 366
$templateRow = $startRow - 1

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $key = $rows[$i][0]
    $val = $rows[$i][1]

    # Copy formatting (cell style) from the last pre-existing data row so the
    # new rows match the rest of the table (style "import"/s="1") instead of
    # getting the workbook default style.
    $ws.Range("A" + $templateRow + ":C" + $templateRow).Copy()
    $ws.Range("A" + $r + ":C" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = "cs"
    $ws.Cells.Item($r, 2).Value = $key
    $ws.Cells.Item($r, 3).Value = $val
}

$excel.ActiveWindow.ScrollRow = 358
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B372").Select()
